$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 773.6842
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 980
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 2940
$ws.Range("M70").Value = -1830
$ws.Range("N70").Value = -3480

$ws.Range("H73").Value = 773.6842
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 980
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 2940
$ws.Range("M73").Value = -1164
$ws.Range("N73").Value = -4812

$ws.Range("H98").Value = 1134.8667
$ws.Range("I98").Value = 1263.3077
$ws.Range("K98").Value = 1263.3077
$ws.Range("M98").Value = 234.6922999999999

$ws.Range("H122").Value = 1134.8667
$ws.Range("I122").Value = 1263.3077
$ws.Range("K122").Value = 3789.9231
$ws.Range("M122").Value = -1339.9231

$ws.Range("H127").Value = 858.7059
$ws.Range("I127").Value = 600.3077
$ws.Range("J127").Value = 1698.5
$ws.Range("K127").Value = 1800.9231
$ws.Range("L127").Value = 5095.5
$ws.Range("M127").Value = 3159.0769
$ws.Range("N127").Value = -15015.5

$ws.Range("H129").Value = 7345.1763
$ws.Range("I129").Value = 421.33334
$ws.Range("K129").Value = 1264.00002
$ws.Range("M129").Value = 3735.99998

$ws.Range("H138").Value = 3476.3076
$ws.Range("I138").Value = 2028.8064
$ws.Range("J138").Value = 4224.183
$ws.Range("K138").Value = 6086.4192
$ws.Range("L138").Value = 12672.549
$ws.Range("M138").Value = -946.4192000000003
$ws.Range("N138").Value = -22952.549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 744.5
$ws.Range("I2").Value = 744.5
$ws.Range("K2").Value = 744.5
$ws.Range("M2").Value = -631.5

$ws.Range("H4").Value = 1455.6666
$ws.Range("I4").Value = 1016.8333
$ws.Range("J4").Value = 2333.3333
$ws.Range("K4").Value = 1016.8333
$ws.Range("L4").Value = 2333.3333
$ws.Range("M4").Value = -900.8333
$ws.Range("N4").Value = -2565.3333

$ws.Range("H9").Value = 17950
$ws.Range("J9").Value = 17950
$ws.Range("L9").Value = 17950
$ws.Range("N9").Value = -18290

$ws.Range("H20").Value = 17950
$ws.Range("J20").Value = 17950
$ws.Range("L20").Value = 17950
$ws.Range("N20").Value = -18490

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = 0

$ws.Range("H32").Value = 23816.535
$ws.Range("I32").Value = 12280.056
$ws.Range("J32").Value = 44582.2
$ws.Range("K32").Value = 12280.056
$ws.Range("L32").Value = 44582.2
$ws.Range("M32").Value = -11993.056
$ws.Range("N32").Value = -45156.2

$ws.Range("H74").Value = 1074.6938
$ws.Range("I74").Value = 868.7381
$ws.Range("J74").Value = 2310.4285
$ws.Range("K74").Value = 868.7381
$ws.Range("L74").Value = 2310.4285
$ws.Range("M74").Value = 5.261899999999969
$ws.Range("N74").Value = -4058.4285

$ws.Range("H77").Value = 1074.6938
$ws.Range("I77").Value = 868.7381
$ws.Range("J77").Value = 2310.4285
$ws.Range("K77").Value = 4343.690500000001
$ws.Range("L77").Value = 11552.1425
$ws.Range("M77").Value = 24.30949999999939
$ws.Range("N77").Value = -20288.1425

$ws.Range("H116").Value = 744.5
$ws.Range("I116").Value = 744.5
$ws.Range("K116").Value = 744.5
$ws.Range("M116").Value = 1549.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 744.5
$ws.Range("I3").Value = 744.5
$ws.Range("K3").Value = 744.5
$ws.Range("M3").Value = -630.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 153.63637
$ws.Range("I22").Value = 144.44444
$ws.Range("J22").Value = 195
$ws.Range("K22").Value = 144.44444
$ws.Range("L22").Value = 195
$ws.Range("M22").Value = 205.55556
$ws.Range("N22").Value = -895

$ws.Range("H31").Value = 2619.4238
$ws.Range("I31").Value = 1187.9615
$ws.Range("J31").Value = 3747.2424
$ws.Range("K31").Value = 1187.9615
$ws.Range("L31").Value = 3747.2424
$ws.Range("M31").Value = -892.9614999999999
$ws.Range("N31").Value = -4337.2424

$ws.Range("H34").Value = 2619.4238
$ws.Range("I34").Value = 1187.9615
$ws.Range("J34").Value = 3747.2424
$ws.Range("K34").Value = 1187.9615
$ws.Range("L34").Value = 3747.2424
$ws.Range("M34").Value = -985.9614999999999
$ws.Range("N34").Value = -4151.2424

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 9998
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 9998
$ws.Range("K75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("M75").Value = 29994
$ws.Range("N75").Value = -31990

$ws.Range("H78").Value = 9998
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 9998
$ws.Range("K78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("M78").Value = 89982
$ws.Range("N78").Value = -99966

$ws.Range("H123").Value = 3270.1155
$ws.Range("I123").Value = 1132
$ws.Range("J123").Value = 4606.4375
$ws.Range("K123").Value = 3396
$ws.Range("L123").Value = 13819.3125
$ws.Range("M123").Value = -946
$ws.Range("N123").Value = -18719.3125

$ws.Range("H131").Value = 921.2222
$ws.Range("J131").Value = 1169.6875
$ws.Range("L131").Value = 3509.0625
$ws.Range("N131").Value = -13589.0625

$ws.Range("H132").Value = 1114.4375
$ws.Range("J132").Value = 1628.5714
$ws.Range("L132").Value = 14657.1426
$ws.Range("N132").Value = -19717.1426

$ws.Range("H134").Value = 842.1667
$ws.Range("I134").Value = 613.25
$ws.Range("K134").Value = 1839.75
$ws.Range("M134").Value = 3230.25

$ws.Range("H139").Value = 5206.1787
$ws.Range("I139").Value = 1935.2142
$ws.Range("J139").Value = 8477.143
$ws.Range("K139").Value = 5805.642599999999
$ws.Range("L139").Value = 25431.429
$ws.Range("M139").Value = -665.6425999999992
$ws.Range("N139").Value = -35711.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1401.9565
$ws.Range("I122").Value = 1492.5
$ws.Range("K122").Value = 4477.5
$ws.Range("M122").Value = -2027.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2066.6667
$ws.Range("I7").Value = 2066.6667
$ws.Range("K7").Value = 2066.6667
$ws.Range("M7").Value = -1954.6667

$ws.Range("H36").Value = 29000
$ws.Range("J36").Value = 29000
$ws.Range("L36").Value = 29000
$ws.Range("N36").Value = -30124

$ws.Range("H40").Value = 2048.6
$ws.Range("I40").Value = 1964.5
$ws.Range("K40").Value = 1964.5
$ws.Range("M40").Value = -1828.5

$ws.Range("H126").Value = 2066.6667
$ws.Range("I126").Value = 2066.6667
$ws.Range("K126").Value = 6200.000100000001
$ws.Range("M126").Value = -3730.000100000001

$ws.Range("H132").Value = 8202222.5
$ws.Range("I132").Value = 9096556
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 27289668
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -27287138
$ws.Range("N132").Value = -17559.0005

$ws.Range("H136").Value = 5678.9062
$ws.Range("I136").Value = 5840.24
$ws.Range("J136").Value = 5102.7144
$ws.Range("K136").Value = 17520.72
$ws.Range("L136").Value = 15308.1432
$ws.Range("M136").Value = -14970.72
$ws.Range("N136").Value = -20408.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1939.2307
$ws.Range("I126").Value = 2296.111
$ws.Range("J126").Value = 1136.25
$ws.Range("K126").Value = 6888.333
$ws.Range("L126").Value = 3408.75
$ws.Range("M126").Value = -4418.333
$ws.Range("N126").Value = -8348.75

$ws.Range("H132").Value = 1951
$ws.Range("I132").Value = 726.6857
$ws.Range("K132").Value = 2180.0571
$ws.Range("M132").Value = 349.9429

$ws.Range("H136").Value = 3411.3333
$ws.Range("I136").Value = 1010.2857
$ws.Range("J136").Value = 11815
$ws.Range("K136").Value = 3030.8571
$ws.Range("L136").Value = 35445
$ws.Range("M136").Value = -480.8571000000002
$ws.Range("N136").Value = -40545
